# Apply the AOCS-sheet hardware-design-choice rework described in the commit
# "Model fully working, all hardware design choices available".
#
# Summary of the change:
#  - AOCS sheet ("name"/"value"/"units"/"description" table): columns C and D
#    (units / description) were swapped so the header now reads
#    name | value | description | units, for every one of the 14 top rows.
#  - "P mass" / "P volume" rows were renamed "P mass w/ HS" / "P volume w/ HS"
#    and got new (lander heat-shield-less) values.
#  - Three brand new rows were appended for the lander hardware choice:
#    P mass lander / P volume lander / P power lander.
#  - Cosmetic: AOCS column A got an explicit width, the AOCS sheet selection
#    moved to C31, and the Astro sheet's column A width/selection changed too.

$wb = $excel.ActiveWorkbook
$wsAocs = $wb.Worksheets.Item("AOCS")
$wsAstro = $wb.Worksheets.Item("Astro")

# --- AOCS: swap the "units" and "description" columns (C <-> D) for the
# header row and the 13 data rows that have both a unit and a description.
for ($r = 1; $r -le 14; $r++) {
    $cVal = $wsAocs.Cells.Item($r, 3).Value2
    $dVal = $wsAocs.Cells.Item($r, 4).Value2
    $wsAocs.Cells.Item($r, 3).Value = $dVal
    $wsAocs.Cells.Item($r, 4).Value = $cVal
}

# --- AOCS: rename the probe mass/volume rows to the "w/ HS" (with heat
# shield) variants and refresh their numbers.
$wsAocs.Range("A18").Value = "P mass w/ HS"
$wsAocs.Range("B18").Value = 3.37
$wsAocs.Range("A19").Value = "P volume w/ HS"
$wsAocs.Range("B19").Value = 0.0145

# --- AOCS: append the new lander rows.
$wsAocs.Range("A20").Value = "P mass lander"
$wsAocs.Range("B20").Value = 34.4
$wsAocs.Range("A21").Value = "P volume lander"
$wsAocs.Range("B21").Value = 0.31759999999999999
$wsAocs.Range("A22").Value = "P power lander"
$wsAocs.Range("B22").Value = 94

# --- Cosmetic view/format touch-ups -----------------------------------

# Astro sheet: widen column A a bit (manual width, autofit marker dropped)
# and leave the cursor on A10.
$wsAstro.Columns.Item(1).ColumnWidth = 17.5
$wsAstro.Range("A10").Select()

# AOCS: give column A an explicit width and finish with the selection on
# C31 (also re-activates the AOCS tab, as in the source file).
$wsAocs.Columns.Item(1).ColumnWidth = 14.666666666666666
$wsAocs.Range("C31").Select()
